# Rename Sheet1 -> erosion
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "erosion"

# Add a new worksheet for soil data, placed after the erosion sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "soil"

# Fill header row for soil sheet
$ws2.Range("B1").Value = "PH值"
$ws2.Range("C1").Value = "全氮含量(%)"
$ws2.Range("D1").Value = "碳酸根离子"
$ws2.Range("E1").Value = "硫酸根离子"
$ws2.Range("F1").Value = "镁离子"
$ws2.Range("G1").Value = "钾离子"
$ws2.Range("H1").Value = "钠离子"
$ws2.Range("I1").Value = "土壤电阻(Ω)"
$ws2.Range("J1").Value = "站点"

# Fill data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 6.75
$ws2.Range("C2").Value = 0.109
$ws2.Range("D2").Value = 0.0126
$ws2.Range("E2").Value = 0.0118
$ws2.Range("F2").Value = 0.0017
$ws2.Range("G2").Value = 0.0004
$ws2.Range("H2").Value = 0.0049
$ws2.Range("I2").Value = 32.9
$ws2.Range("J2").Value = "沈阳站"

# Match the source selection/active sheet state
$ws2.Range("A1:J2").Select() | Out-Null
$ws2.Activate() | Out-Null
